# Apply updated cryptocurrency price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '65.192.48'
$ws.Range('E2').Value = '  -2.24%  '

# Row 3
$ws.Range('D3').Value = '3.475.39'
$ws.Range('E3').Value = '  -1.34%  '

# Row 4
$ws.Range('E4').Value = '  +0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '589.38'
$ws.Range('E5').Value = '  -2.83%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '137.16'
$ws.Range('E6').Value = '  -4.69%  '

# Row 7
$ws.Range('D7').Value = '3.476.08'
$ws.Range('E7').Value = '  -1.30%  '

# Row 8
$ws.Range('E8').Value = '  +0.05%  '

# Row 9
$ws.Range('E9').Value = '  -3.12%  '

# Row 10
$ws.Range('E10').Value = '  -6.11%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.15'
$ws.Range('E11').Value = '  -6.99%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.383'
$ws.Range('E12').Value = '  -5.19%  '

# Row 13
$ws.Range('D13').Value = '4.066.27'
$ws.Range('E13').Value = '  -1.13%  '

# Row 14
$ws.Range('E14').Value = '  -6.92%  '

# Row 15
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '26.60'
$ws.Range('E15').Value = '  -7.33%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.474.51'
$ws.Range('E16').Value = '  -1.27%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.115'
$ws.Range('E17').Value = '  -1.23%  '

# Row 18
$ws.Range('D18').Value = '65.155.27'
$ws.Range('E18').Value = '  -2.01%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.70'
$ws.Range('E19').Value = '  -9.63%  '

# Row 20
$ws.Range('E20').Value = '  -5.89%  '

# Row 21
$ws.Range('E21').Value = '  -4.97%  '

# Row 23
$ws.Range('E23').Value = '  -5.89%  '

# Row 24
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.78'
$ws.Range('E24').Value = '  +1.00%  '

# Row 25
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '72.62'
$ws.Range('E25').Value = '  -5.52%  '

# Row 26
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.03%  '

# Row 27
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '3.619.06'
$ws.Range('E27').Value = '  -1.34%  '

# Row 28
$ws.Range('E28').Value = '  -4.17%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.05%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.30'
$ws.Range('E30').Value = '  -6.54%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.17'
$ws.Range('E31').Value = '  -8.49%  '

# Row 32
$ws.Range('E32').Value = '  -9.76%  '

# Row 33
$ws.Range('D33').Value = '3.496.10'
$ws.Range('E33').Value = '  -0.92%  '

# Row 34
$ws.Range('E34').Value = '  -0.06%  '

# Row 35
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '23.14'
$ws.Range('E35').Value = '  -4.39%  '

# Row 36
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.143'
$ws.Range('E36').Value = '  -7.55%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '170.96'
$ws.Range('E37').Value = '  -1.25%  '

# Row 38
$ws.Range('E38').Value = '  -10.13%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.82'
$ws.Range('E39').Value = '  -9.58%  '

# Row 40
$ws.Range('E40').Value = '  -9.25%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.71'
$ws.Range('E41').Value = '  -9.22%  '

# Row 42
$ws.Range('E42').Value = '  -3.21%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.811'

# Row 44
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.10%  '

# Row 45
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '42.43'
$ws.Range('E45').Value = '  -6.92%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '25.14'
$ws.Range('E46').Value = '  +9.72%  '

# Row 47
$ws.Range('E47').Value = '  -12.51%  '

# Row 48
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.62'
$ws.Range('E48').Value = '  -8.16%  '

# Row 49
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.15'
$ws.Range('E49').Value = '  +3.09%  '

# Row 50
$ws.Range('E50').Value = '  -5.08%  '

# Row 51
$ws.Range('D51').Value = '2.219.54'
$ws.Range('E51').Value = '  -3.64%  '
